$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.495.74"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +6.78%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.483.60"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.58%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "490.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.71%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.43"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +15.87%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("E8").Value = "  +8.13%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.503.26"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.84%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.80"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +10.46%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0980"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.82%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.333"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +8.07%  "

$ws.Range("E13").Value = "  +1.91%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.915.94"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.52%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "56.413.51"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.38%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +10.13%  "

$ws.Range("E17").Value = "  +6.17%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.499.10"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.47%  "

$ws.Range("E19").Value = "  +11.35%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +10.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "320.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.77%  "

$ws.Range("E22").Value = "  +0.40%  "

$ws.Range("E23").Value = "  +10.87%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "58.77"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.54%  "

$ws.Range("E25").Value = "  +8.68%  "

$ws.Range("E26").Value = "  -0.71%  "

$ws.Range("E27").Value = "  +9.48%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.593.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.88%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.63"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +9.63%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0793"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +12.29%  "

$ws.Range("E31").Value = "  +0.25%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "149.46"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.46%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.25"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.79%  "

$ws.Range("E34").Value = "  +7.86%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.24"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.48%  "

$ws.Range("E36").Value = "  +10.19%  "

$ws.Range("E37").Value = "  +8.03%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.861"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +10.11%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "34.25"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.93%  "

$ws.Range("E40").Value = "  +9.09%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0562"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.74%  "

$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.612"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.50%  "

$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.993"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.13%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.34"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +9.96%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.83"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +17.50%  "

$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0924"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.23%  "

$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "259.67"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +20.11%  "

$ws.Range("E48").Value = "  +7.15%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.21"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.09%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.66"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.88%  "

$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.893.62"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.80%  "
